# Refresh the crypto price/volume figures as described in the commit diff.
#
# Each data row (2-51) holds one coin; normally only its Price (column D)
# and Volume(1h) (column E) text change. Rows 46 and 47 additionally swap
# their Coin/Link content (Maker <-> Celestia) in addition to getting new
# Price/Volume figures.
#
# Price values that look like plain numbers (e.g. "1.00", "0.999") are
# written with a leading apostrophe so Excel stores them as literal text
# (preserving e.g. trailing zeros) instead of silently coercing them into
# numeric cells. Immediately afterwards the cell Style is reset back to
# "Normal" so no stray number-format/quote-prefix style sticks to the
# cell, keeping it identical in style to the original (unstyled) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.362.34'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '2.277.56'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = "'310.30"
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("D6").Value = "'104.20"
$ws.Range("E6").Value = '  +3.18%  '
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("D10").Value = "'38.71"
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").Value = "'0.0898"
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = "'8.21"
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = "'0.973"
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("D15").Value = "'15.02"
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").Value = '2.616.64'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = '2.271.58'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '42.566.83'
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").Value = "'7.23"
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("E20").Value = '  -0.74%  '
$ws.Range("D21").Value = "'13.38"
$ws.Range("E21").Value = '  +5.92%  '
$ws.Range("D22").Value = "'72.80"
$ws.Range("E23").Value = '  -3.69%  '
$ws.Range("D24").Value = "'262.22"
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("D27").Value = "'10.63"
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = "'6.93"
$ws.Range("E29").Value = '  +15.25%  '
$ws.Range("D30").Value = "'22.17"
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").Value = "'35.65"
$ws.Range("E31").Value = '  -5.04%  '
$ws.Range("D32").Value = "'164.30"
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = "'0.0851"
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("E35").Value = '  +2.35%  '
$ws.Range("E36").Value = '  -2.49%  '
$ws.Range("D37").Value = "'4.49"
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("D38").Value = "'0.0346"
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("D40").Value = "'2.72"
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("D41").Value = "'1.54"
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("D42").Value = "'98.09"
$ws.Range("E42").Value = '  +8.09%  '
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").Value = "'68.39"
$ws.Range("E44").Value = '  +0.66%  '
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").Value = "'11.89"
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '1.713.20'
$ws.Range("E47").Value = '  +6.22%  '
$ws.Range("D48").Value = "'109.91"
$ws.Range("E48").Value = '  -4.07%  '
$ws.Range("D49").Value = "'75.96"
$ws.Range("E49").Value = '  -3.50%  '
$ws.Range("D50").Value = "'5.14"
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").Value = "'8.60"
$ws.Range("E51").Value = '  -3.52%  '

# Cells above that were forced to text with a quote-prefix: strip the
# resulting style override back to the default so the cell matches the
# original (no explicit style) formatting.
$textForcedCells = @(
    "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12",
    "D13", "D14", "D15", "D19", "D21", "D22", "D24", "D27",
    "D29", "D30", "D31", "D32", "D33", "D37", "D38", "D40",
    "D41", "D42", "D43", "D44", "D46", "D48", "D49", "D50",
    "D51"
)
foreach ($ref in $textForcedCells) {
    $ws.Range($ref).Style = "Normal"
}

